# InfraFair control inputs - correction of EU example (description text fixes)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value  = "The percentage of demand responsibility for the cost of the assets"
$ws.Range("D4").Value  = "The percentage of generation responsibility for the cost of the assets"
$ws.Range("D8").Value  = "The number of hours each snapshot represents, the total should be 8760 hours (one year)"
$ws.Range("D9").Value  = "1 to allocate the full cost, 2 to allocate only the cost of the used capacity, 3 to allocate full cost if the assets are classified as 'Exist' and the cost of the used capacity if they are classified as 'Planned', 4 is to allocate the cost based on the utilization threshold, if the asset is utilized more than the threshold, allocate the full cost, otherwise, allocate the cost of the used capacity"
$ws.Range("D10").Value = "If the ratio between the used asset capacity and the asset rated capacity is equal or above this percentage, the asset cost will be fully allocated, otherwise, only the cost of the used capacity will be allocated.This will be used only if 'Cost Allocation Option' is set to 4"
$ws.Range("D11").Value = "To determine what to do with the cost of unused capacity in case the 'Cost Allocation Option' is not set to 1. 0 to do nothing, 1 to allocate it equally among agents who use the asset, 2 to allocate it equally among all agents of the country(ies) owning the asset, 3 to allocate it equally among all agents"
$ws.Range("D12").Value = "The percentage of demand responsibility to the socialized cost of the assets. Only used when the 'Cost Allocation Option' is not set to 1 and 'Cost of Unused Capacity' is not set to 0"
$ws.Range("D13").Value = "The percentage of generators responsibility for the socialized cost of the assets. Only used when 'Cost Allocation Option' is not set to 1 and 'Cost of Unused Capacity' is not set to 0"

# Reflect the final cursor / scroll position left by the edit (last cell touched was D10)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D10").Select()
